# SysSettings.xlsx edit: "Change transmission network loss - all regions are
# divided into the smaller regions"
#
# The "Constants" sheet had a single region column ("AllRegions") that is
# being split into two regions ("DKE" and "DKW"). This is implemented as a
# native column insert before column F, which:
#   - shifts every column F..CR right by one (F..CR -> G..CS), automatically
#     re-pointing every formula in the workbook that referenced the moved
#     cells (on this sheet and on others, e.g. the LOG sheet's ADDRESS()
#     formulas);
#   - leaves a blank column F that we populate with the new "DKW" region's
#     data (mirroring the existing region columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# --- 1. Insert the new region column before column F -----------------------
$ws.Columns("F:F").Insert()

# --- 2. Header row (row 6): rename AllRegions -> DKE, add DKW --------------
# New F6 should carry the same formatting as the other region-header cells
# (e.g. G6, which used to be F6 before the insert).
$ws.Range("G6").Copy() | Out-Null
$ws.Range("F6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("E6").Value = "DKE"
$ws.Range("F6").Value = "DKW"

# --- 3. Row 7 (price year row): chain F7=E7, G7=F7, ... K7=J7 --------------
$ws.Range("F7:K7").NumberFormat = "0"
$ws.Range("F7").Formula = "=E7"
$ws.Range("G7").Formula = "=F7"
$ws.Range("H7").Formula = "=G7"
$ws.Range("I7").Formula = "=H7"
$ws.Range("J7").Formula = "=I7"
$ws.Range("K7").Formula = "=J7"

# --- 4. Rows 8-17 (conversion factor table): chain formulas, new col F gets
#        the same number format ("0.000") as the rest of the row -----------
$ws.Range("F8:F17").NumberFormat = "0.000"
for ($r = 8; $r -le 17; $r++) {
    $ws.Cells.Item($r, 6).Formula  = "=" + $ws.Cells.Item($r, 5).Address($false, $false)
    $ws.Cells.Item($r, 7).Formula  = "=" + $ws.Cells.Item($r, 6).Address($false, $false)
    $ws.Cells.Item($r, 8).Formula  = "=" + $ws.Cells.Item($r, 7).Address($false, $false)
    $ws.Cells.Item($r, 9).Formula  = "=" + $ws.Cells.Item($r, 8).Address($false, $false)
    $ws.Cells.Item($r, 10).Formula = "=" + $ws.Cells.Item($r, 9).Address($false, $false)
    $ws.Cells.Item($r, 11).Formula = "=" + $ws.Cells.Item($r, 10).Address($false, $false)
}

# --- 5. Row 18 (discount-rate-ish row of plain values): fill 0.1 across ----
$ws.Range("F18:K18").NumberFormat = "General"
$ws.Range("F18:K18").Value = 0.1

# --- 6. Row 19: new region columns get their own format (0.0) + value 1 ---
$ws.Range("F19:K19").NumberFormat = "0.0"
$ws.Range("F19:K19").Value = 1

$excel.CutCopyMode = 0
